$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# E1: "Categoria" -> "Email"
$ws.Range("E1").Value = "Email"
# F1 (new column): "Category"
$ws.Range("F1").Value = "Category"

# --- New column F for existing data rows 2 and 3 ---
$ws.Range("F2").Value = "Alimentaria"
$ws.Range("F3").Value = "Alimentaria"

# --- New rows 4, 5, 6 ---
# Build them by copying row 3 (which already has the D/E "numeric-looking"
# text values stored as shared strings) so the text cell type is preserved
# without Excel re-interpreting "2255789"/"45454" as numbers and without
# introducing any new cell styles.
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A4:F4").PasteSpecial(-4104) | Out-Null

$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A5:F5").PasteSpecial(-4104) | Out-Null

$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A6:F6").PasteSpecial(-4104) | Out-Null

# Now adjust the cells that differ from the copied row 3 template.
$ws.Range("A5").Value = "CAMPERO"
$ws.Range("A6").Value = "CAMPERO"
$ws.Range("E6").Value = "45454sdsdsd"
